$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from the placeholder "Tets" to "Sheet1"
$ws.Name = "Sheet1"

# Replace the template header row with the real alumni data row
$ws.Range("A1").Value = "subhan assiddik"

$ws.Range("B1").ClearFormats()
$ws.Range("B1").Value = 201721021

$ws.Range("C1").Value = "subhanassiddik@gmail.com"

$ws.Range("D1").ClearFormats()
$ws.Range("D1").Value = 85256199655

# Drop the mailto hyperlink's explicit display text so it just tracks the cell text
foreach ($hl in $ws.Hyperlinks) {
    $hl.TextToDisplay = ""
}

# The template had two extra blank rows below the header - remove them
$ws.Rows("2:3").Delete()

# E1 only ever held formatting (no value) in the template - clear it out
$ws.Range("E1").Clear()

# Move the visible selection to E1, matching the saved view state
$ws.Range("E1").Select() | Out-Null
